$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 209, shifting rows 209:254 down to 210:255.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new data record.
$ws.Cells.Item(209, 1).Value = 3
$ws.Cells.Item(209, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 45005
$ws.Cells.Item(209, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(209, 5).Value = 5
$ws.Cells.Item(209, 6).Value = 100112052
$ws.Cells.Item(209, 7).Value = "Albahaca"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 93
$ws.Cells.Item(209, 11).Value = 4500
$ws.Cells.Item(209, 12).Value = 5000
$ws.Cells.Item(209, 13).Value = 4796
$ws.Cells.Item(209, 14).Value = "`$/docena de matas"
$ws.Cells.Item(209, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(209, 16).Value = 799
$ws.Cells.Item(209, 17).Value = 6
$ws.Cells.Item(209, 18).Value = "Hortaliza"
